# Scheduled-runner refresh of market-board derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on the
# per-job Leve sheets. Values below are the new pulls; only the cells
# that actually moved are touched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1180.591
$ws.Range("I19").Value = 811.75
$ws.Range("J19").Value = 1623.2
$ws.Range("K19").Value = 811.75
$ws.Range("L19").Value = 1623.2
$ws.Range("M19").Value = -636.75
$ws.Range("N19").Value = -1973.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1668.1428
$ws.Range("I28").Value = 1277.3334
$ws.Range("J28").Value = 2840.5715
$ws.Range("K28").Value = 1277.3334
$ws.Range("L28").Value = 2840.5715
$ws.Range("M28").Value = -792.3334
$ws.Range("N28").Value = -3810.5715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 726.16
$ws.Range("I38").Value = 41.333332
$ws.Range("J38").Value = 2487.1428
$ws.Range("K38").Value = 123.999996
$ws.Range("L38").Value = 7461.428400000001
$ws.Range("M38").Value = 248.000004
$ws.Range("N38").Value = -8205.428400000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 565.3913
$ws.Range("I41").Value = 637.2143
$ws.Range("J41").Value = 453.66666
$ws.Range("K41").Value = 637.2143
$ws.Range("L41").Value = 453.66666
$ws.Range("M41").Value = -197.2143
$ws.Range("N41").Value = -1333.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1458.0667
$ws.Range("I53").Value = 4725.25
$ws.Range("J53").Value = 270
$ws.Range("K53").Value = 4725.25
$ws.Range("L53").Value = 270
$ws.Range("M53").Value = -4088.25
$ws.Range("N53").Value = -1544

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3190
$ws.Range("I76").Value = 3190
$ws.Range("K76").Value = 3190
$ws.Range("M76").Value = -2875

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3190
$ws.Range("I79").Value = 3190
$ws.Range("K79").Value = 3190
$ws.Range("M79").Value = -2098

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2000
$ws.Range("J111").Value = 2000
$ws.Range("L111").Value = 6000
$ws.Range("N111").Value = -12134

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1287.5714
$ws.Range("I129").Value = 900
$ws.Range("J129").Value = 1442.6
$ws.Range("K129").Value = 2700
$ws.Range("L129").Value = 4327.799999999999
$ws.Range("M129").Value = 2300
$ws.Range("N129").Value = -14327.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6349
$ws.Range("I132").Value = 2387.0588
$ws.Range("J132").Value = 15970.857
$ws.Range("K132").Value = 7161.176399999999
$ws.Range("L132").Value = 47912.571
$ws.Range("M132").Value = -4631.176399999999
$ws.Range("N132").Value = -52972.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12837.554
$ws.Range("I32").Value = 7169.381
$ws.Range("J32").Value = 23188.13
$ws.Range("K32").Value = 7169.381
$ws.Range("L32").Value = 23188.13
$ws.Range("M32").Value = -6882.381
$ws.Range("N32").Value = -23762.13

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 731674.25
$ws.Range("I134").Value = 1114929.2
$ws.Range("J134").Value = 5506.737
$ws.Range("K134").Value = 3344787.6
$ws.Range("L134").Value = 16520.211
$ws.Range("M134").Value = -3342252.6
$ws.Range("N134").Value = -21590.211

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 611.6667
$ws.Range("I16").Value = 534
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 534
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -247
$ws.Range("N16").Value = -1574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 611.6667
$ws.Range("I113").Value = 534
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 534
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1636
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 947544.5
$ws.Range("I113").Value = 3030899.2
$ws.Range("J113").Value = 565
$ws.Range("K113").Value = 9092697.600000001
$ws.Range("L113").Value = 1695
$ws.Range("M113").Value = -9090527.600000001
$ws.Range("N113").Value = -6035

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 925.47
$ws.Range("J131").Value = 956.883
$ws.Range("L131").Value = 2870.649
$ws.Range("N131").Value = -12950.649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5428.3335
$ws.Range("I70").Value = 5220
$ws.Range("J70").Value = 5688.75
$ws.Range("K70").Value = 5220
$ws.Range("L70").Value = 5688.75
$ws.Range("M70").Value = -4950
$ws.Range("N70").Value = -6228.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5428.3335
$ws.Range("I73").Value = 5220
$ws.Range("J73").Value = 5688.75
$ws.Range("K73").Value = 5220
$ws.Range("L73").Value = 5688.75
$ws.Range("M73").Value = -4284
$ws.Range("N73").Value = -7560.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 994.375
$ws.Range("I122").Value = 751.2
$ws.Range("J122").Value = 1399.6666
$ws.Range("K122").Value = 2253.6
$ws.Range("L122").Value = 4198.9998
$ws.Range("M122").Value = 196.3999999999996
$ws.Range("N122").Value = -9098.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1739.25
$ws.Range("I7").Value = 1719
$ws.Range("J7").Value = 1800
$ws.Range("K7").Value = 1719
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = -1607
$ws.Range("N7").Value = -2024

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1002.51514
$ws.Range("I46").Value = 907.7692
$ws.Range("J46").Value = 1064.1
$ws.Range("K46").Value = 907.7692
$ws.Range("L46").Value = 1064.1
$ws.Range("M46").Value = -719.7692
$ws.Range("N46").Value = -1440.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3230.7693
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1739.25
$ws.Range("I126").Value = 1719
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 5157
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2687
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 3007.5
$ws.Range("I54").Value = 3826.6667
$ws.Range("J54").Value = 550
$ws.Range("K54").Value = 3826.6667
$ws.Range("L54").Value = 550
$ws.Range("M54").Value = -3306.6667
$ws.Range("N54").Value = -1590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 8867.583000000001
$ws.Range("I107").Value = 12901.375
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 38704.125
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -36784.125
$ws.Range("N107").Value = -6240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 29994.5
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
